# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.014.16'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.633.18'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'214.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('D9').Value = "'0.0620"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.21%  '
$ws.Range('D10').Value = "'18.38"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.46%  '
$ws.Range('D11').Value = "'0.0792"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').Value = '1.630.71'
$ws.Range('E13').Value = '  -2.06%  '
$ws.Range('D14').Value = "'4.18"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = "'0.525"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.54%  '
$ws.Range('D16').Value = '25.992.13'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').Value = '0.0₃0741'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').Value = "'61.41"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.20%  '
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').Value = "'190.96"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').Value = "'9.64"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('D23').Value = "'6.07"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.68%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = "'1.80"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('D26').Value = "'144.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = "'6.77"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').Value = "'15.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.40%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('D31').Value = "'0.0481"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('D32').Value = "'3.14"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('E33').Value = '  -5.51%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'1.49"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'2.41"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').Value = '1.130.44'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').Value = "'0.861"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.15%  '
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').Value = "'0.518"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').Value = "'98.38"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').Value = "'0.775"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('D43').Value = '1.769.40'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').Value = "'54.86"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  -3.57%  '
